$wb = $excel.ActiveWorkbook

# --- Sheet "Técnicos" (1st sheet): update two names ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B4").Value = "Zé 2"
$ws1.Range("B5").Value = "Manel 3"

# --- Sheet "Projetos" (2nd sheet): update two project names ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ANI1"
$ws2.Range("B3").Value = "Ani2"

# --- Update selection on Técnicos sheet (no longer the active tab) ---
$ws1.Range("B6").Select()

# --- Make Projetos the active sheet/tab and set its selection ---
$ws2.Activate()
$ws2.Range("B4").Select()
